$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "445÷2=222, 1" "674÷8=84, 2"
Replace-Text "779÷7=111, 2" "701÷5=140, 1"
Replace-Text "132÷9=14, 6" "946÷7=135, 1"
Replace-Text "798÷7=114, 0" "194÷2=97, 0"
Replace-Text "640÷3=213, 1" "830÷3=276, 2"
Replace-Text "297÷4=74, 1" "495÷6=82, 3"
Replace-Text "523÷7=74, 5" "437÷5=87, 2"
Replace-Text "363÷3=121, 0" "522÷8=65, 2"
Replace-Text "303÷7=43, 2" "934÷8=116, 6"
Replace-Text "591÷4=147, 3" "509÷8=63, 5"
Replace-Text "457÷6=76, 1" "297÷8=37, 1"
Replace-Text "116÷2=58, 0" "614÷5=122, 4"
Replace-Text "696÷5=139, 1" "234÷4=58, 2"
Replace-Text "432÷7=61, 5" "447÷7=63, 6"
Replace-Text "883÷4=220, 3" "511÷2=255, 1"
Replace-Text "430÷2=215, 0" "600÷6=100, 0"
Replace-Text "870÷7=124, 2" "246÷4=61, 2"
Replace-Text "944÷6=157, 2" "239÷4=59, 3"
Replace-Text "449÷9=49, 8" "128÷2=64, 0"
Replace-Text "439÷9=48, 7" "712÷2=356, 0"
Replace-Text "252÷5=50, 2" "356÷8=44, 4"
Replace-Text "165÷2=82, 1" "848÷7=121, 1"
Replace-Text "415÷2=207, 1" "393÷2=196, 1"
Replace-Text "305÷3=101, 2" "936÷3=312, 0"
Replace-Text "759÷7=108, 3" "507÷4=126, 3"

Write-Output "Done"
